$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "839÷9="; New = "376÷9=" },
    @{ Row = 1;  Col = 2; Old = "831÷9="; New = "546÷4=" },
    @{ Row = 1;  Col = 3; Old = "222÷4="; New = "390÷8=" },
    @{ Row = 1;  Col = 4; Old = "857÷9="; New = "516÷3=" },
    @{ Row = 1;  Col = 5; Old = "104÷9="; New = "793÷7=" },

    @{ Row = 5;  Col = 1; Old = "522÷6="; New = "834÷4=" },
    @{ Row = 5;  Col = 2; Old = "527÷9="; New = "608÷3=" },
    @{ Row = 5;  Col = 3; Old = "859÷5="; New = "427÷3=" },
    @{ Row = 5;  Col = 4; Old = "920÷3="; New = "980÷5=" },
    @{ Row = 5;  Col = 5; Old = "484÷6="; New = "136÷6=" },

    @{ Row = 9;  Col = 1; Old = "295÷3="; New = "180÷3=" },
    @{ Row = 9;  Col = 2; Old = "135÷4="; New = "901÷3=" },
    @{ Row = 9;  Col = 3; Old = "926÷4="; New = "895÷7=" },
    @{ Row = 9;  Col = 4; Old = "288÷7="; New = "168÷7=" },
    @{ Row = 9;  Col = 5; Old = "960÷5="; New = "688÷9=" },

    @{ Row = 13; Col = 1; Old = "138÷7="; New = "979÷6=" },
    @{ Row = 13; Col = 2; Old = "589÷3="; New = "382÷5=" },
    @{ Row = 13; Col = 3; Old = "348÷7="; New = "546÷8=" },
    @{ Row = 13; Col = 4; Old = "312÷4="; New = "903÷8=" },
    @{ Row = 13; Col = 5; Old = "841÷5="; New = "200÷8=" },

    @{ Row = 17; Col = 1; Old = "217÷5="; New = "106÷8=" },
    @{ Row = 17; Col = 2; Old = "841÷5="; New = "662÷8=" },
    @{ Row = 17; Col = 3; Old = "962÷6="; New = "491÷2=" },
    @{ Row = 17; Col = 4; Old = "291÷8="; New = "489÷2=" },
    @{ Row = 17; Col = 5; Old = "685÷6="; New = "328÷2=" }
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $rng = $cell.Range
    # Replace:=1 (wdReplaceOne) scopes the replacement to the single match
    # found within this cell's Range, avoiding accidental document-wide
    # replacement of duplicate values (e.g. "841÷5=" occurs twice).
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 1)
}
